$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 299.3
$ws.Range("B3").Value = -188.4
$ws.Range("C3").Value = 28.4
$ws.Range("C4").Value = 3
$ws.Range("C20").Value = -54.9
$ws.Range("C21").Value = -109.7
$ws.Range("C22").Value = -36.4
